# Adds 5 new worksheets (addProductToBasket, trackOurOrder, searchProduct,
# currysStoreFinder, otherCurrysWebsiteNavigation) with their test data,
# mirroring the commit "Added 5 test cases and corresponding classes."

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# Sheet: addProductToBasket
# ---------------------------------------------------------------------------
$ws3 = Add-SheetAtEnd "addProductToBasket"

$ws3.Cells.Item(1,1).Value = "browser"
$ws3.Cells.Item(1,2).Value = "menu name"
$ws3.Cells.Item(1,3).Value = "product name"
$ws3.Cells.Item(1,4).Value = "product title"
$ws3.Cells.Item(1,5).Value = "item brand"
$ws3.Cells.Item(1,6).Value = "item name"
$ws3.Cells.Item(1,7).Value = "item added to basket success message"

$ws3.Cells.Item(2,1).Value = "chrome"
$ws3.Cells.Item(2,2).Value = "Appliances"
$ws3.Cells.Item(2,3).Value = "Washing machines"
$ws3.Cells.Item(2,4).Value = "Washing machines"
$ws3.Cells.Item(2,5).Value = "HOTPOINT"
$ws3.Cells.Item(2,6).Value = "NSWM 1043C GG UK N 10 kg 1400 Spin Washing Machine - Graphite"
$ws3.Cells.Item(2,7).Value = "This item has been added to your basket"

$ws3.Cells.Item(3,1).Value = "firefox"
$ws3.Cells.Item(3,2).Value = "Appliances"
$ws3.Cells.Item(3,3).Value = "Kettles"
$ws3.Cells.Item(3,4).Value = "Kettles"
$ws3.Cells.Item(3,5).Value = "ESSENTIALS"
$ws3.Cells.Item(3,6).Value = "C17JKW17 Jug Kettle - White"
$ws3.Cells.Item(3,7).Value = "This item has been added to your basket"

$ws3.Cells.Item(4,1).Value = "edge"
$ws3.Cells.Item(4,2).Value = "TV & Audio"
$ws3.Cells.Item(4,3).Value = "Soundbars"
$ws3.Cells.Item(4,4).Value = "Sound bars"
$ws3.Cells.Item(4,5).Value = "LOGIK"
$ws3.Cells.Item(4,6).Value = "L32SBIN16A 2.1 Sound Bar"
$ws3.Cells.Item(4,7).Value = "This item has been added to your basket"

$ws3.Columns.Item(2).ColumnWidth = 9.983072916666666
$ws3.Columns.Item(3).ColumnWidth = 15.619791666666666
$ws3.Columns.Item(4).ColumnWidth = 15.619791666666666
$ws3.Columns.Item(5).ColumnWidth = 9.529947916666666
$ws3.Columns.Item(6).ColumnWidth = 57.072916666666664
$ws3.Columns.Item(7).ColumnWidth = 34.166666666666664

$ws3.Range("D11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: trackOurOrder
# ---------------------------------------------------------------------------
$ws4 = Add-SheetAtEnd "trackOurOrder"

$ws4.Cells.Item(1,1).Value = "browser"
$ws4.Cells.Item(1,2).Value = "job or repair ref"
$ws4.Cells.Item(1,3).Value = "search by type"
$ws4.Cells.Item(1,4).Value = "search by data"
$ws4.Cells.Item(1,5).Value = "country"

# The author typed rows 2-3 with the range already formatted as Text, so
# numeric-looking strings (e.g. the phone number "78128190") are preserved
# verbatim instead of being auto-coerced to numbers -- format *before*
# writing the values.
$ws4.Range("A2:E3").NumberFormat = "@"

$ws4.Cells.Item(2,1).Value = "edge"
$ws4.Cells.Item(2,2).Value = 8122332
$ws4.Cells.Item(2,3).Value = "Contact Number"
$ws4.Cells.Item(2,4).Value = "78128190"
$ws4.Cells.Item(2,5).Value = "Republic of Ireland"

$ws4.Cells.Item(3,1).Value = "chrome"
$ws4.Cells.Item(3,2).Value = "A689912"
$ws4.Cells.Item(3,3).Value = "Postcode"
$ws4.Cells.Item(3,4).Value = "UER123"
$ws4.Cells.Item(3,5).Value = "United Kingdom"

$ws4.Columns.Item(2).ColumnWidth = 13.256510416666666
$ws4.Columns.Item(3).ColumnWidth = 13.709635416666666
$ws4.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws4.Columns.Item(5).ColumnWidth = 15.709635416666666

$ws4.Range("H15").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: searchProduct
# ---------------------------------------------------------------------------
$ws5 = Add-SheetAtEnd "searchProduct"

$ws5.Cells.Item(1,1).Value = "browser"
$ws5.Cells.Item(1,2).Value = "search data"
$ws5.Cells.Item(1,3).Value = "page Title contains"

$ws5.Cells.Item(2,1).Value = "edge"
$ws5.Cells.Item(2,2).Value = "grinder"
$ws5.Cells.Item(2,3).Value = "grinder"

$ws5.Columns.Item(2).ColumnWidth = 9.619791666666666
$ws5.Columns.Item(3).ColumnWidth = 14.893229166666666

$ws5.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: currysStoreFinder
# ---------------------------------------------------------------------------
$ws6 = Add-SheetAtEnd "currysStoreFinder"

$ws6.Cells.Item(1,1).Value = "browser"
$ws6.Cells.Item(1,2).Value = "town / postcode"
$ws6.Cells.Item(1,3).Value = "expected store"

$ws6.Cells.Item(2,1).Value = "chrome"
$ws6.Cells.Item(2,2).Value = "penzance"
$ws6.Cells.Item(2,3).Value = "Currys, Penzance"

$ws6.Cells.Item(3,1).Value = "edge"
$ws6.Cells.Item(3,2).Value = "UB34FF"
$ws6.Cells.Item(3,3).Value = "Currys, Hayes"

$ws6.Columns.Item(2).ColumnWidth = 13.709635416666666
$ws6.Columns.Item(3).ColumnWidth = 17.893229166666668

$ws6.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: otherCurrysWebsiteNavigation
# ---------------------------------------------------------------------------
$ws7 = Add-SheetAtEnd "otherCurrysWebsiteNavigation"

$ws7.Cells.Item(1,1).Value = "browser"
$ws7.Cells.Item(1,2).Value = "website"
$ws7.Cells.Item(1,3).Value = "expected page title"
$ws7.Cells.Item(1,4).Value = "expected URL"

$ws7.Cells.Item(2,1).Value = "edge"
$ws7.Cells.Item(2,2).Value = "Currys Business"
$ws7.Cells.Item(2,3).Value = "Welcome - Currys Business"
$ws7.Cells.Item(2,4).Value = "https://business.currys.co.uk/"

$ws7.Cells.Item(3,1).Value = "chrome"
$ws7.Cells.Item(3,2).Value = "Currys Ireland"
$ws7.Cells.Item(3,3).Value = "Currys PC World | Laptops, TVs, Washing Machines, Cookers, Smartphones & Lots More"
$ws7.Cells.Item(3,4).Value = "https://www.currys.ie/ieen/index.html"

$ws7.Cells.Item(4,1).Value = "mozilla"
$ws7.Cells.Item(4,2).Value = "Partmaster"
$ws7.Cells.Item(4,3).Value = "Appliance Parts & Electrical Accessories for Home Appliances | Currys Partmaster.co.uk"
$ws7.Cells.Item(4,4).Value = "https://www.partmaster.co.uk/"

# Turn column D (2..4) into real hyperlinks pointing at the URL that is
# already displayed as the cell's own text -- keep the existing cell value
# as the on-screen text by not passing TextToDisplay.
$ws7.Hyperlinks.Add($ws7.Range("D2"), "https://business.currys.co.uk/") | Out-Null
$ws7.Hyperlinks.Add($ws7.Range("D3"), "https://www.currys.ie/ieen/index.html") | Out-Null
$ws7.Hyperlinks.Add($ws7.Range("D4"), "https://www.partmaster.co.uk/") | Out-Null

$ws7.Columns.Item(2).ColumnWidth = 13.072916666666666
$ws7.Columns.Item(3).ColumnWidth = 73.98307291666667
$ws7.Columns.Item(4).ColumnWidth = 33.072916666666664

$ws7.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore the selection on the pre-existing "pageNavigation" sheet, then make
# the newly-added last sheet ("otherCurrysWebsiteNavigation") the active tab
# again, matching the saved workbook view.
# ---------------------------------------------------------------------------
$wsNav = $wb.Worksheets.Item("pageNavigation")
$wsNav.Activate()
$wsNav.Range("C11").Select() | Out-Null

$ws7.Activate()
